$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Insert a new row at 406 (pushes old rows 406..522 down to 407..523) ---
$ws.Rows("406").Insert()

# Copy the formatting (styles/borders) from the row above (405) into the
# newly inserted blank row 406 so it matches the rest of the table rows.
$ws.Range("A405:K405").Copy()
$ws.Range("A406:K406").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Grow Table1 so it covers the new row (was A8:K522, now A8:K523) ---
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K523"))

# Restore the calculated-column formula on the last (formerly-last) row,
# which Excel mangled into a #VALUE! error when the row temporarily fell
# outside the table during the insert/resize above.
$ws.Range("G523").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),""."",Table1[[#This Row],[EARNED]])"
$ws.Range("G523").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- Fill in the new leave-card entries ---

# Row 403 (09/01/2023): SL earned 1.25
$ws.Range("C403").Value = 1.25

# Row 404 (10/01/2023): SL(1-0-0) taken, SL earned 1.25, 1 day absent w/ pay, dated 10/13/2023
$ws.Range("B404").Value = "SL(1-0-0)"
$ws.Range("C404").Value = 1.25
$ws.Range("H404").Value = 1
$ws.Range("K404").Value = 45212

# Row 405 (11/01/2023): SL(1-0-0) taken, SL earned 1.25, 1 day absent w/ pay, dated 11/20/2023
$ws.Range("B405").Value = "SL(1-0-0)"
$ws.Range("C405").Value = 1.25
$ws.Range("H405").Value = 1
$ws.Range("K405").Value = 45250

# Row 406 (new): Forced leave FL(3-0-0), 3 days absence w/o pay, remarks 12/27-29/2023
$ws.Range("B406").Value = "FL(3-0-0)"
$ws.Range("D406").Value = 3
$ws.Range("G406").ClearContents()
$ws.Range("K406").Value = "12/27-29/2023"

# K404/K405/K406 need the existing "date" style (same one already used e.g.
# by K400) instead of the plain style copied from row 405 -- copy it over.
$ws.Range("K400").Copy()
$ws.Range("K404:K406").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("K406").Value = "12/27-29/2023"

# --- Recalculate everything (refreshes G-column / E9 / I9 totals) ---
$excel.CalculateFullRebuild()
